$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column holds price strings; prefix with an apostrophe so Excel
# keeps them as text instead of inferring a number (preserves exact
# formatting such as trailing zeros / thousand-dot grouping).

$ws.Range('D2').Value = "'51.773.79"
$ws.Range('E2').Value = '  +5.25%  '
$ws.Range('D3').Value = "'2.763.13"
$ws.Range('E3').Value = '  +5.43%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'116.39"
$ws.Range('E5').Value = '  +3.92%  '
$ws.Range('D6').Value = "'332.76"
$ws.Range('E6').Value = '  +3.10%  '
$ws.Range('E7').Value = '  +2.46%  '
$ws.Range('D8').Value = "'0.999"
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = "'0.577"
$ws.Range('E9').Value = '  +6.46%  '
$ws.Range('D10').Value = "'41.88"
$ws.Range('E10').Value = '  +5.36%  '
$ws.Range('D11').Value = "'0.0865"
$ws.Range('E11').Value = '  +6.88%  '
$ws.Range('D12').Value = "'20.19"
$ws.Range('E12').Value = '  +2.32%  '
$ws.Range('E13').Value = '  +2.05%  '
$ws.Range('E14').Value = '  +5.42%  '
$ws.Range('D15').Value = "'3.194.75"
$ws.Range('E15').Value = '  +5.67%  '
$ws.Range('D16').Value = "'2.753.24"
$ws.Range('E16').Value = '  +4.51%  '
$ws.Range('D17').Value = "'0.888"
$ws.Range('E17').Value = '  +3.76%  '
$ws.Range('D18').Value = "'51.678.88"
$ws.Range('E18').Value = '  +5.27%  '
$ws.Range('D19').Value = "'3.21"
$ws.Range('E19').Value = '  +6.34%  '
$ws.Range('E20').Value = '  +4.58%  '
$ws.Range('D21').Value = "'6.87"
$ws.Range('E21').Value = '  +2.90%  '
$ws.Range('D22').Value = "'0.0₃0981"
$ws.Range('E22').Value = '  +4.11%  '
$ws.Range('D23').Value = "'278.41"
$ws.Range('E24').Value = '  +1.78%  '
$ws.Range('E25').Value = '  +5.28%  '
$ws.Range('E26').Value = '  +2.62%  '
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('D28').Value = "'10.19"
$ws.Range('E28').Value = '  -1.02%  '
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('E30').Value = '  +1.92%  '
$ws.Range('D31').Value = "'35.16"
$ws.Range('E31').Value = '  +1.11%  '
$ws.Range('D32').Value = "'50.00"
$ws.Range('E32').Value = '  +0.91%  '
$ws.Range('E33').Value = '  +1.87%  '
$ws.Range('D34').Value = "'0.0824"
$ws.Range('E34').Value = '  +2.01%  '
$ws.Range('D35').Value = "'0.999"
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('D36').Value = "'19.00"
$ws.Range('E36').Value = '  -0.07%  '
$ws.Range('D37').Value = "'5.01"
$ws.Range('E37').Value = '  +1.68%  '
$ws.Range('E38').Value = '  +2.67%  '
$ws.Range('E39').Value = '  +3.47%  '
$ws.Range('D40').Value = "'0.0352"
$ws.Range('E40').Value = '  +10.10%  '
$ws.Range('D41').Value = "'126.99"
$ws.Range('E41').Value = '  +0.16%  '
$ws.Range('D42').Value = "'23.14"
$ws.Range('E42').Value = '  +5.23%  '
$ws.Range('E43').Value = '  +3.28%  '
$ws.Range('E44').Value = '  +7.55%  '
$ws.Range('D45').Value = "'2.45"
$ws.Range('E45').Value = '  +13.65%  '
$ws.Range('D46').Value = "'2.090.74"
$ws.Range('E46').Value = '  +1.69%  '
$ws.Range('D47').Value = "'3.32"
$ws.Range('E47').Value = '  +3.31%  '
$ws.Range('E48').Value = '  +4.91%  '
$ws.Range('E49').Value = '  +6.88%  '
$ws.Range('D50').Value = "'9.04"
$ws.Range('E50').Value = '  +1.54%  '
$ws.Range('D51').Value = "'60.06"
$ws.Range('E51').Value = '  +1.93%  '
